# Auto-generated edit script: update crypto price/volume snapshot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.011.25'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").Value = '2.473.83'
$ws.Range("E3").Value = '  +0.82%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '560.57'
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.91'
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.505'
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +0.27%  '
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.333'
$ws.Range("E11").Value = '  -2.72%  '
$ws.Range("E12").Value = '  +1.32%  '
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("D14").Value = '68.912.91'
$ws.Range("E14").Value = '  +0.77%  '
$ws.Range("E15").Value = '  -1.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.62'
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("D17").Value = '2.473.29'
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.71'
$ws.Range("E18").Value = '  -2.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '335.86'
$ws.Range("E19").Value = '  -2.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.96'
$ws.Range("E20").Value = '  -3.11%  '
$ws.Range("E21").Value = '  -0.59%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("B23").Value = 'SuiNetwork'
$ws.Range("C23").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.88'
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.77'
$ws.Range("E24").Value = '  -1.97%  '
$ws.Range("D25").Value = '2.601.19'
$ws.Range("E25").Value = '  +1.16%  '
$ws.Range("E26").Value = '  -2.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -5.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.18'
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("D29").Value = '0.0₃0818'
$ws.Range("E29").Value = '  -2.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.20'
$ws.Range("E30").Value = '  -1.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '430.21'
$ws.Range("E32").Value = '  -1.47%  '
$ws.Range("E33").Value = '  -4.17%  '
$ws.Range("E34").Value = '  -4.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '159.03'
$ws.Range("E35").Value = '  +1.42%  '
$ws.Range("E36").Value = '  +0.18%  '
$ws.Range("E39").Value = '  -0.66%  '
$ws.Range("E40").Value = '  -2.31%  '
$ws.Range("E41").Value = '  -1.77%  '
$ws.Range("E42").Value = '  -4.34%  '
$ws.Range("E43").Value = '  -2.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.07'
$ws.Range("E44").Value = '  -1.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '131.71'
$ws.Range("E45").Value = '  -2.58%  '
$ws.Range("E46").Value = '  -0.63%  '
$ws.Range("E47").Value = '  -0.92%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.484'
$ws.Range("E48").Value = '  -1.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.563'
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0912'
$ws.Range("E50").Value = '  -0.27%  '
$ws.Range("E51").Value = '  +0.38%  '
